$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.07830000000001
$ws.Range("D4").Value = -7.051200000000005
$ws.Range("D5").Value = -8.272299999999996
$ws.Range("A7").Value = -21.56290000000001
$ws.Range("D8").Value = -8.272299999999996
$ws.Range("A16").Value = -20.23339999999998
$ws.Range("D16").Value = -8.249999999999998
